$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 169
$ws.Cells.Item(4, 6).Value = 246
$ws.Cells.Item(5, 6).Value = 700848
$ws.Cells.Item(6, 6).Value = 1937
$ws.Cells.Item(7, 6).Value = 1655
$ws.Cells.Item(8, 6).Value = 720
$ws.Cells.Item(9, 6).Value = 572
$ws.Cells.Item(10, 6).Value = 1450
$ws.Cells.Item(11, 6).Value = 1393
$ws.Cells.Item(12, 6).Value = 1218
$ws.Cells.Item(13, 6).Value = 2892
$ws.Cells.Item(14, 6).Value = 1691
$ws.Cells.Item(15, 6).Value = 1102
$ws.Cells.Item(16, 6).Value = 1681
$ws.Cells.Item(17, 6).Value = 48
$ws.Cells.Item(18, 6).Value = 8
$ws.Cells.Item(19, 6).Value = 9
$ws.Cells.Item(20, 6).Value = 578
$ws.Cells.Item(21, 6).Value = 1148
$ws.Cells.Item(22, 6).Value = 1650
$ws.Cells.Item(23, 6).Value = 1650
$ws.Cells.Item(25, 6).Value = 59
$ws.Cells.Item(26, 6).Value = 152
$ws.Cells.Item(27, 6).Value = 2136
$ws.Cells.Item(28, 6).Value = 1633
$ws.Cells.Item(30, 6).Value = 4276
$ws.Cells.Item(35, 6).Value = 219
$ws.Cells.Item(36, 6).Value = 312
$ws.Cells.Item(37, 6).Value = 68
$ws.Cells.Item(38, 6).Value = 2809
$ws.Cells.Item(39, 6).Value = 234
$ws.Cells.Item(40, 6).Value = 1091
$ws.Cells.Item(41, 6).Value = 3333
$ws.Cells.Item(42, 6).Value = 1079
$ws.Cells.Item(43, 6).Value = 42
$ws.Cells.Item(45, 6).Value = 202
$ws.Cells.Item(47, 6).Value = 25
$ws.Cells.Item(48, 6).Value = 56
$ws.Cells.Item(49, 6).Value = 722
$ws.Cells.Item(50, 6).Value = 29
$ws.Cells.Item(51, 6).Value = 47

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(2, 6).Value = 68
$ws.Cells.Item(3, 6).Value = 71
$ws.Cells.Item(4, 6).Value = 4
$ws.Cells.Item(5, 6).Value = 162
$ws.Cells.Item(7, 6).Value = 148025
$ws.Cells.Item(8, 6).Value = 148025
$ws.Cells.Item(10, 6).Value = 64
$ws.Cells.Item(12, 6).Value = 23
$ws.Cells.Item(13, 6).Value = 111
$ws.Cells.Item(14, 6).Value = 242
$ws.Cells.Item(16, 6).Value = 9
$ws.Cells.Item(22, 6).Value = 94
$ws.Cells.Item(23, 6).Value = 945
$ws.Cells.Item(26, 6).Value = 12
$ws.Cells.Item(28, 6).Value = 71
$ws.Cells.Item(29, 6).Value = 458
$ws.Cells.Item(30, 6).Value = 303
$ws.Cells.Item(31, 6).Value = 13
$ws.Cells.Item(32, 6).Value = 119
$ws.Cells.Item(33, 6).Value = 119
$ws.Cells.Item(36, 6).Value = 253
$ws.Cells.Item(37, 6).Value = 115
$ws.Cells.Item(39, 6).Value = 196
$ws.Cells.Item(40, 6).Value = 94
$ws.Cells.Item(44, 6).Value = 1

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 121
$ws.Cells.Item(4, 6).Value = 3185
$ws.Cells.Item(5, 6).Value = 295
$ws.Cells.Item(7, 6).Value = 878
$ws.Cells.Item(8, 6).Value = 1313
$ws.Cells.Item(9, 6).Value = 687
$ws.Cells.Item(10, 6).Value = 242
$ws.Cells.Item(11, 6).Value = 2311

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 295
$ws.Cells.Item(3, 6).Value = 878
$ws.Cells.Item(4, 6).Value = 1313
$ws.Cells.Item(5, 6).Value = 687
$ws.Cells.Item(6, 6).Value = 169
$ws.Cells.Item(7, 6).Value = 246
$ws.Cells.Item(8, 6).Value = 242
$ws.Cells.Item(9, 6).Value = 2311
$ws.Cells.Item(10, 6).Value = 700848
$ws.Cells.Item(11, 6).Value = 162
$ws.Cells.Item(12, 6).Value = 524
$ws.Cells.Item(13, 6).Value = 1937
$ws.Cells.Item(14, 6).Value = 148025
$ws.Cells.Item(15, 6).Value = 1655
$ws.Cells.Item(16, 6).Value = 720
$ws.Cells.Item(17, 6).Value = 572
$ws.Cells.Item(18, 6).Value = 1393
$ws.Cells.Item(19, 6).Value = 1218
$ws.Cells.Item(20, 6).Value = 2892
$ws.Cells.Item(21, 6).Value = 1691
$ws.Cells.Item(22, 6).Value = 30
$ws.Cells.Item(23, 6).Value = 1681
$ws.Cells.Item(24, 6).Value = 578
$ws.Cells.Item(25, 6).Value = 233
$ws.Cells.Item(26, 6).Value = 1148
$ws.Cells.Item(27, 6).Value = 1650
$ws.Cells.Item(28, 6).Value = 1650
$ws.Cells.Item(29, 6).Value = 1168
$ws.Cells.Item(30, 6).Value = 96
$ws.Cells.Item(31, 6).Value = 152
$ws.Cells.Item(32, 6).Value = 2136
$ws.Cells.Item(33, 6).Value = 1633
$ws.Cells.Item(35, 6).Value = 4276
$ws.Cells.Item(36, 6).Value = 36
$ws.Cells.Item(37, 6).Value = 1193
$ws.Cells.Item(38, 6).Value = 103
$ws.Cells.Item(39, 6).Value = 458
$ws.Cells.Item(40, 6).Value = 219
$ws.Cells.Item(41, 6).Value = 303
$ws.Cells.Item(43, 6).Value = 312
$ws.Cells.Item(44, 6).Value = 2809
$ws.Cells.Item(45, 6).Value = 234
$ws.Cells.Item(46, 6).Value = 1091
$ws.Cells.Item(47, 6).Value = 3333
$ws.Cells.Item(48, 6).Value = 253
$ws.Cells.Item(49, 6).Value = 1079
$ws.Cells.Item(51, 6).Value = 202
$ws.Cells.Item(52, 6).Value = 699
$ws.Cells.Item(53, 6).Value = 56
$ws.Cells.Item(54, 6).Value = 722
